$wb = $excel.ActiveWorkbook

# Best-effort: remember the last on-screen window position/size (cosmetic,
# mirrors the author's local Excel window geometry).
$win = $excel.ActiveWindow
$win.Left = 6500
$win.Top = 5020
$win.Width = 19640
$win.Height = 14540

$ws = $wb.Worksheets.Item("materials")

# Add four new material/resource columns (F:I) with their headers and
# starting values, matching the new columns introduced for the
# "citizen on house upgrade" fix.
$ws.Range("F1").Value = "INT_trainingFigure"
$ws.Range("G1").Value = "INT_bowTarget"
$ws.Range("H1").Value = "INT_saddle"
$ws.Range("I1").Value = "INT_ironPart"

$ws.Range("F2").Value = 1000
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000

# Move the selection/active cell on the "materials" sheet to D7.
$ws.Activate() | Out-Null
$ws.Range("D7").Select() | Out-Null
